$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2226776683031646
$ws.Range("C2").Value = 0.5618053799771385
$ws.Range("D2").Value = 0.6246702391763954
$ws.Range("E2").Value = 0.7903608284678558
$ws.Range("F2").Value = 0.7667232613880856
$ws.Range("G2").Value = 46

$ws.Range("B3").Value = 0.293438081967545
$ws.Range("C3").Value = 0.6695240432108801
$ws.Range("D3").Value = 0.9182272854563196
$ws.Range("E3").Value = 0.9582417677477431
$ws.Range("F3").Value = 0.9225146215623552
$ws.Range("G3").Value = 45

$ws.Range("B4").Value = 0.182530046921793
$ws.Range("C4").Value = 0.5838921936173209
$ws.Range("D4").Value = 0.6581603071260018
$ws.Range("E4").Value = 0.8112707976539041
$ws.Range("F4").Value = 0.7996088567086325
$ws.Range("G4").Value = 44

$ws.Range("B5").Value = 0.327844975786004
$ws.Range("C5").Value = 0.7216035043336653
$ws.Range("D5").Value = 1.079816604652694
$ws.Range("E5").Value = 1.03914224466754
$ws.Range("F5").Value = 0.9977399924889347
$ws.Range("G5").Value = 43

$ws.Range("B6").Value = 0.2304610288088512
$ws.Range("C6").Value = 0.5964001741841488
$ws.Range("D6").Value = 0.7431306070296015
$ws.Range("E6").Value = 0.8620502346322989
$ws.Range("F6").Value = 0.840742550595554
$ws.Range("G6").Value = 42

$ws.Range("B7").Value = 0.360861710175336
$ws.Range("C7").Value = 0.7403235333647489
$ws.Range("D7").Value = 1.140448443186462
$ws.Range("E7").Value = 1.067917807317802
$ws.Range("F7").Value = 1.0175868272775
$ws.Range("G7").Value = 41

$ws.Range("B8").Value = 0.2216292893798094
$ws.Range("C8").Value = 0.5674340590213579
$ws.Range("D8").Value = 0.6895196896333891
$ws.Range("E8").Value = 0.8303732230951267
$ws.Range("F8").Value = 0.8104447324343937
$ws.Range("G8").Value = 40

$ws.Range("B9").Value = 0.4304460383242469
$ws.Range("C9").Value = 0.7770904355848363
$ws.Range("D9").Value = 1.255534736424988
$ws.Range("E9").Value = 1.120506464249532
$ws.Range("F9").Value = 1.048053168048187
$ws.Range("G9").Value = 39

$ws.Range("B10").Value = 0.2885592243077175
$ws.Range("C10").Value = 0.5810449934012875
$ws.Range("D10").Value = 0.7354929410298064
$ws.Range("E10").Value = 0.8576088508345786
$ws.Range("F10").Value = 0.818446246706525
$ws.Range("G10").Value = 38

$ws.Range("B11").Value = 0.4217362328470671
$ws.Range("C11").Value = 0.7352401290682791
$ws.Range("D11").Value = 1.136174966075426
$ws.Range("E11").Value = 1.065915083895254
$ws.Range("F11").Value = 0.9924380765910319
$ws.Range("G11").Value = 37
